# Update "Forecast Comparison" sheet with the correct forecast output:
#  - insert a new "Week_Start_Date" column after "Week"
#  - populate it with the week's start date (as text, e.g. "2025-01-05")
#  - renumber the "Week" labels from zero-padded (W01..W09) to plain (W1..W9)
#  - the trailing "is_holiday_week" column becomes a boolean (FALSE) column

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new blank column at B; everything from the old ASIN column onward
# (ASIN, MyForecast, Amazon Mean/P70/P80/P90 Forecast, Product Title,
# is_holiday_week) shifts one column to the right automatically.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Make column B plain text so the date-like strings aren't coerced into
# Excel date serial numbers.
$ws.Columns.Item(2).NumberFormat = "@"

$weekStartDates = @{
    2  = "2025-01-05"
    3  = "2025-01-12"
    4  = "2025-01-19"
    5  = "2025-01-26"
    6  = "2025-02-02"
    7  = "2025-02-09"
    8  = "2025-02-16"
    9  = "2025-02-23"
    10 = "2025-03-02"
    11 = "2025-03-09"
    12 = "2025-03-16"
    13 = "2025-03-23"
    14 = "2025-03-30"
    15 = "2025-04-06"
    16 = "2025-04-13"
    17 = "2025-04-20"
}

foreach ($row in 2..17) {
    # Week_Start_Date column.
    $ws.Cells.Item($row, 2).Value = $weekStartDates[$row]

    # Week labels drop their leading zero: "W01".."W09" -> "W1".."W9".
    $ws.Cells.Item($row, 1).Value = "W" + ($row - 1)

    # is_holiday_week (now column J) becomes a real boolean FALSE instead
    # of the numeric 0 it used to hold.
    $ws.Cells.Item($row, 10).Value = $false
}
